# Apply month-append changes to the indirect-expenses workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B1: report month moved from Apr-2024 (45383) to May-2024 (45413)
$ws.Range("B1").Value = 45413

# Updated figures for existing rows
$ws.Range("B2").Value  = 368463.4
$ws.Range("B3").Value  = 132151.5
$ws.Range("B7").Value  = 432
$ws.Range("B9").Value  = 22149
$ws.Range("B10").Value = 12300
$ws.Range("B13").Value = 9002.9699999999993
$ws.Range("B20").Value = 6420
$ws.Range("B21").Value = 13283.2
$ws.Range("B22").Value = 19032
$ws.Range("B23").Value = 5500
$ws.Range("B29").Value = -1.55
$ws.Range("B33").Value = 10509
$ws.Range("B35").Value = 6559.62
$ws.Range("B36").Value = 5410
$ws.Range("B38").Value = 21180
$ws.Range("B39").Value = 64166.66
$ws.Range("B40").Value = 810580.95
$ws.Range("B41").Value = 3610.26
$ws.Range("B42").Value = 172490.69
$ws.Range("B43").Value = 105226
$ws.Range("B45").Value = 399870
$ws.Range("B47").Value = 10384
$ws.Range("B51").Value = 50771.509999999995
$ws.Range("B55").Value = 22.96
$ws.Range("B56").Value = 22632.55
$ws.Range("B60").Value = 16896
$ws.Range("B64").Value = 1229815.8599999999
$ws.Range("B66").Value = 368463
$ws.Range("B67").Value = 50772
$ws.Range("B68").Value = 419235
$ws.Range("B69").Value = 575971
$ws.Range("B71").Value = 10384
$ws.Range("B72").Value = 224226
$ws.Range("B73").Value = 810581
$ws.Range("B74").Value = 1229816

# Rows that now gain a value where previously empty
$ws.Range("B24").Value = 1500
$ws.Range("B27").Value = 1500
$ws.Range("B30").Value = 6000
$ws.Range("B54").Value = 10000
$ws.Range("B58").Value = 1220

# Rows whose value is removed entirely (cell becomes empty again)
$ws.Range("B12").ClearContents()
$ws.Range("B14").ClearContents()
$ws.Range("B25").ClearContents()
$ws.Range("B26").ClearContents()
